$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand new row at position 83; this shifts the former rows 83-125
# down to 84-126 (dimension grows from A1:R125 to A1:R126).
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new weekly price record.
$ws.Range("A83").Value = 8
$ws.Range("B83").Value = "Terminal La Palmera de La Serena"
$ws.Range("C83").Value = "Coquimbo"
$ws.Range("D83").Value = (Get-Date -Year 2023 -Month 3 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E83").Value = 4
$ws.Range("F83").Value = 100112030
$ws.Range("G83").Value = "Poroto granado"
$ws.Range("H83").Value = "Sin especificar"
$ws.Range("I83").Value = "Primera"
$ws.Range("J83").Value = 450
$ws.Range("K83").Value = 34500
$ws.Range("L83").Value = 35000
$ws.Range("M83").Value = 34750
$ws.Range("N83").Value = "$/malla 25 kilos"
$ws.Range("O83").Value = "Provincia del Elquí"
$ws.Range("P83").Value = 1390
$ws.Range("Q83").Value = 25
$ws.Range("R83").Value = "Hortaliza"
